# Auto-generated: applies 2025-05-31 data update to violent-crime-full-year workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 2611
$ws.Range("L3").Value = 2616
$ws.Range("H4").Value = 1756
$ws.Range("I4").Value = 1835
$ws.Range("L4").Value = 700
$ws.Range("L5").Value = 154
$ws.Range("L6").Value = 2345
$ws.Range("H7").Value = 26070
$ws.Range("I7").Value = 26303
$ws.Range("L7").Value = 8426

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L7").Value = 275
$ws.Range("L8").Value = 536
$ws.Range("L10").Value = 55
$ws.Range("L11").Value = 147
$ws.Range("L12").Value = 20
$ws.Range("L13").Value = 12
$ws.Range("L19").Value = 236
$ws.Range("L20").Value = 217
$ws.Range("L24").Value = 19
$ws.Range("L25").Value = 50
$ws.Range("L29").Value = 447
$ws.Range("L30").Value = 40
$ws.Range("L33").Value = 372
$ws.Range("L37").Value = 313
$ws.Range("L42").Value = 276
$ws.Range("L46").Value = 19
$ws.Range("L48").Value = 116
$ws.Range("L50").Value = 45
$ws.Range("L51").Value = 97
$ws.Range("L52").Value = 166
$ws.Range("L54").Value = 168
$ws.Range("H63").Value = 305
$ws.Range("I63").Value = 259
$ws.Range("L63").Value = 23
$ws.Range("L64").Value = 55
$ws.Range("L67").Value = 313
$ws.Range("L76").Value = 105
$ws.Range("L78").Value = 111
$ws.Range("L79").Value = 228
$ws.Range("L84").Value = 87
$ws.Range("L85").Value = 437
$ws.Range("L86").Value = 60
$ws.Range("L89").Value = 110
$ws.Range("L91").Value = 120
$ws.Range("L93").Value = 43
$ws.Range("L98").Value = 57
$ws.Range("L99").Value = 141
$ws.Range("H101").Value = 26070
$ws.Range("I101").Value = 26303
$ws.Range("L101").Value = 8426

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L2").Value = 82
$ws.Range("L3").Value = 86
$ws.Range("L6").Value = 77
$ws.Range("L7").Value = 275

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L2").Value = 54
$ws.Range("L6").Value = 36
$ws.Range("L7").Value = 147

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("L2").Value = 37
$ws.Range("L3").Value = 28
$ws.Range("L7").Value = 110

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L3").Value = 177
$ws.Range("L7").Value = 437

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L2").Value = 55
$ws.Range("L7").Value = 166

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L3").Value = 176
$ws.Range("L6").Value = 148
$ws.Range("L7").Value = 536

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L2").Value = 103
$ws.Range("L5").Value = 7
$ws.Range("L7").Value = 372

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L2").Value = 94
$ws.Range("L6").Value = 98
$ws.Range("L7").Value = 313

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L3").Value = 60
$ws.Range("L4").Value = 12
$ws.Range("L5").Value = 2
$ws.Range("L6").Value = 33
$ws.Range("L7").Value = 141

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("L6").Value = 19
$ws.Range("L7").Value = 40

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L2").Value = 94
$ws.Range("L3").Value = 109
$ws.Range("L7").Value = 313

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("L3").Value = 35
$ws.Range("L7").Value = 87

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L2").Value = 37
$ws.Range("L3").Value = 32
$ws.Range("L7").Value = 168

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 143
$ws.Range("L4").Value = 18
$ws.Range("L7").Value = 447

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("L3").Value = 24
$ws.Range("L7").Value = 116

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L2").Value = 78
$ws.Range("L6").Value = 75
$ws.Range("L7").Value = 236

$ws = $wb.Worksheets.Item("River North")
$ws.Range("L2").Value = 20
$ws.Range("L7").Value = 105

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L2").Value = 75
$ws.Range("L7").Value = 276

$ws = $wb.Worksheets.Item("Boystown")
$ws.Range("L3").Value = 5
$ws.Range("L6").Value = 12

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("L2").Value = 26
$ws.Range("L3").Value = 13
$ws.Range("L7").Value = 55

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("L6").Value = 33
$ws.Range("L7").Value = 111

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("L2").Value = 9
$ws.Range("L7").Value = 19

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("L2").Value = 4
$ws.Range("L7").Value = 19

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("L2").Value = 48
$ws.Range("L7").Value = 120

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L3").Value = 84
$ws.Range("L7").Value = 228

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("L2").Value = 19
$ws.Range("L7").Value = 55

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L2").Value = 68
$ws.Range("L7").Value = 217

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("L2").Value = 16
$ws.Range("L7").Value = 43

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("L3").Value = 25
$ws.Range("L7").Value = 50

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("L6").Value = 30
$ws.Range("L7").Value = 57

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("L3").Value = 12
$ws.Range("L7").Value = 45

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("L4").Value = 34
$ws.Range("L7").Value = 60

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("L3").Value = 28
$ws.Range("L7").Value = 97

$ws = $wb.Worksheets.Item("Beverly")
$ws.Range("L2").Value = 5
$ws.Range("L3").Value = 7
$ws.Range("L7").Value = 20

Write-Host "Applied 2025-05-31 crime data updates"